$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 178; this shifts existing rows 178-259 down to 179-260,
# which reproduces the "after" state for all of those rows (confirmed against
# the diff: after-row(N+1) == before-row(N)) and pushes the last existing row
# (259) down to become the new row 260.
$ws.Rows.Item(178).Insert()

# Populate the newly inserted row 178 with its new data. Columns A, B, C, E,
# F, G, H and R keep the same boilerplate values used throughout this sheet
# for "Feria Lagunitas de Puerto Montt" / "Ají" / "Inferno" entries.
$ws.Cells.Item(178, 1).Value = 4
$ws.Cells.Item(178, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(178, 3).Value = "Los Lagos"
$ws.Cells.Item(178, 4).Value = 44704
$ws.Cells.Item(178, 5).Value = 10
$ws.Cells.Item(178, 6).Value = 100112021
$ws.Cells.Item(178, 7).Value = "Ají"
$ws.Cells.Item(178, 8).Value = "Inferno"
$ws.Cells.Item(178, 9).Value = "Primera"
$ws.Cells.Item(178, 10).Value = 30
$ws.Cells.Item(178, 11).Value = 29000
$ws.Cells.Item(178, 12).Value = 29000
$ws.Cells.Item(178, 13).Value = 29000
$ws.Cells.Item(178, 14).Value = "`$/caja 12 kilos"
$ws.Cells.Item(178, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(178, 16).Value = 2417
$ws.Cells.Item(178, 17).Value = 12
$ws.Cells.Item(178, 18).Value = "Hortaliza"
